# Apply the weekly update: insert two new data rows for Espárragos
# (Femacal de La Calera) corresponding to the new week's prices.
#
# Row inserted at sheet row 6  -> date 2022-10-13 (serial 44847)
# Row inserted at sheet row 14 -> date 2022-10-14 (serial 44848)
# All rows below each insertion point shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($Row, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R) {
    $ws.Cells.Item($Row, 1).Value2 = $A
    $ws.Cells.Item($Row, 2).Value2 = $B
    $ws.Cells.Item($Row, 3).Value2 = $C
    $ws.Cells.Item($Row, 4).Value2 = $D
    $ws.Cells.Item($Row, 5).Value2 = $E
    $ws.Cells.Item($Row, 6).Value2 = $F
    $ws.Cells.Item($Row, 7).Value2 = $G
    $ws.Cells.Item($Row, 8).Value2 = $H
    $ws.Cells.Item($Row, 9).Value2 = $I
    $ws.Cells.Item($Row, 10).Value2 = $J
    $ws.Cells.Item($Row, 11).Value2 = $K
    $ws.Cells.Item($Row, 12).Value2 = $L
    $ws.Cells.Item($Row, 13).Value2 = $M
    $ws.Cells.Item($Row, 14).Value2 = $N
    $ws.Cells.Item($Row, 15).Value2 = $O
    $ws.Cells.Item($Row, 16).Value2 = $P
    $ws.Cells.Item($Row, 17).Value2 = $Q
    $ws.Cells.Item($Row, 18).Value2 = $R
}

# --- Insert new row at position 6 (2022-10-13, Primera) ---
$ws.Rows.Item(6).Insert()
Set-DataRow 6 3 "Femacal de La Calera" "Coquimbo" 44847 5 300000000 "Espárragos" "Verde" "Primera" 1110 1400 1500 1450 "`$/kilo" "Provincia de Quillota" 1450 1 "Hortaliza"

# --- Insert new row at position 14 (2022-10-14, Primera) ---
$ws.Rows.Item(14).Insert()
Set-DataRow 14 3 "Femacal de La Calera" "Coquimbo" 44848 5 300000000 "Espárragos" "Verde" "Primera" 1750 1400 1500 1449 "`$/kilo" "Provincia de Quillota" 1449 1 "Hortaliza"

$wb.Save()
